$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update test data values on row 2
$ws.Range("D2").Value = "pruebasregistro48"
$ws.Range("E2").Value = "1234"
$ws.Range("F2").Value = "4321"

# Update numeroDocumento value
$ws.Range("B2").Value = 700100

# Update the selected cell/active selection in the sheet view
$ws.Range("G7").Select()
